# A new weekly price record was added to the data set. It is inserted as a
# new row 116 (pushing the existing rows 116-162 down to 117-163, which is
# why the sheet's dimension grows from A1:R162 to A1:R163).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 116, shifting rows 116-162 down to 117-163.
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new record's data.
$ws.Range("A116").Value = 3
$ws.Range("B116").Value = "Femacal de La Calera"
$ws.Range("C116").Value = "Coquimbo"
$ws.Range("D116").Value = 44636
$ws.Range("E116").Value = 5
$ws.Range("F116").Value = 100112030
$ws.Range("G116").Value = "Poroto granado"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 45
$ws.Range("K116").Value = 23000
$ws.Range("L116").Value = 24000
$ws.Range("M116").Value = 23444
$ws.Range("N116").Value = "$/malla 25 kilos"
$ws.Range("O116").Value = "Provincia de Talca"
$ws.Range("P116").Value = 938
$ws.Range("Q116").Value = 25
$ws.Range("R116").Value = "Hortaliza"
